$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# The "Ready for handoff" status text is shared across sheets (shared string);
# updating it on every sheet that shows it keeps behavior consistent whether
# or not the backing engine dedupes the string.
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime moved forward (new handback report generated).
$zhcn.Range("K2").Value = "2016-08-20 18:58:14"
$dede.Range("K2").Value = "2016-08-20 18:58:20"

# Error Detail cleared now that the handback file matches the latest version.
$zhcn.Range("P2").Value = ""
$dede.Range("P2").Value = ""

# Column width adjustments (target widths rounded to nearest representable
# Excel column width granularity, matching the COM ColumnWidth contract).
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333332

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333332
